# Update cryptocurrency price/volume data (Tue May  9 04:34:07 UTC 2023 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceCell = $ws.Range("D2")
$priceCell.NumberFormat = "@"
$priceCell.Value = "27.643.68"
$priceCell.Style = "Normal"
$ws.Range("E2").Value = "  -2.18%  "

$priceCell = $ws.Range("D3")
$priceCell.NumberFormat = "@"
$priceCell.Value = "1.844.99"
$priceCell.Style = "Normal"
$ws.Range("E3").Value = "  -1.14%  "

$priceCell = $ws.Range("D4")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.9993"
$priceCell.Style = "Normal"
$ws.Range("E4").Value = "  -0.34%  "

$priceCell = $ws.Range("D5")
$priceCell.NumberFormat = "@"
$priceCell.Value = "314.67"
$priceCell.Style = "Normal"
$ws.Range("E5").Value = "  -1.43%  "

$priceCell = $ws.Range("D6")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.9992"
$priceCell.Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "

$priceCell = $ws.Range("D7")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.4249"
$priceCell.Style = "Normal"
$ws.Range("E7").Value = "  -2.77%  "

$priceCell = $ws.Range("D8")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.3655"
$priceCell.Style = "Normal"
$ws.Range("E8").Value = "  -1.28%  "

$priceCell = $ws.Range("D9")
$priceCell.NumberFormat = "@"
$priceCell.Value = "45.71"
$priceCell.Style = "Normal"
$ws.Range("E9").Value = "  +1.61%  "

$priceCell = $ws.Range("D10")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.07281"
$priceCell.Style = "Normal"
$ws.Range("E10").Value = "  -3.16%  "

$priceCell = $ws.Range("D11")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.9000"
$priceCell.Style = "Normal"
$ws.Range("E11").Value = "  -4.30%  "

$priceCell = $ws.Range("D12")
$priceCell.NumberFormat = "@"
$priceCell.Value = "20.68"
$priceCell.Style = "Normal"
$ws.Range("E12").Value = "  -3.25%  "

$priceCell = $ws.Range("D13")
$priceCell.NumberFormat = "@"
$priceCell.Value = "1.800.65"
$priceCell.Style = "Normal"
$ws.Range("E13").Value = "  -4.95%  "

$priceCell = $ws.Range("D14")
$priceCell.NumberFormat = "@"
$priceCell.Value = "5.388"
$priceCell.Style = "Normal"
$ws.Range("E14").Value = "  -1.21%  "

$priceCell = $ws.Range("D15")
$priceCell.NumberFormat = "@"
$priceCell.Value = "6.569"
$priceCell.Style = "Normal"
$ws.Range("E15").Value = "  -2.30%  "

$priceCell = $ws.Range("D16")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.06869"
$priceCell.Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "

$priceCell = $ws.Range("D17")
$priceCell.NumberFormat = "@"
$priceCell.Value = "1.001"
$priceCell.Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "

$priceCell = $ws.Range("D18")
$priceCell.NumberFormat = "@"
$priceCell.Value = "78.34"
$priceCell.Style = "Normal"
$ws.Range("E18").Value = "  -4.93%  "

$priceCell = $ws.Range("D19")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.000008871"
$priceCell.Style = "Normal"
$ws.Range("E19").Value = "  -2.60%  "

$priceCell = $ws.Range("D20")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.9991"
$priceCell.Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "

$priceCell = $ws.Range("D21")
$priceCell.NumberFormat = "@"
$priceCell.Value = "15.56"
$priceCell.Style = "Normal"
$ws.Range("E21").Value = "  -2.73%  "

$priceCell = $ws.Range("D22")
$priceCell.NumberFormat = "@"
$priceCell.Value = "27.631.48"
$priceCell.Style = "Normal"
$ws.Range("E22").Value = "  -2.22%  "

$priceCell = $ws.Range("D23")
$priceCell.NumberFormat = "@"
$priceCell.Value = "4.974"
$priceCell.Style = "Normal"
$ws.Range("E23").Value = "  -3.01%  "

$priceCell = $ws.Range("D24")
$priceCell.NumberFormat = "@"
$priceCell.Value = "10.61"
$priceCell.Style = "Normal"
$ws.Range("E24").Value = "  -1.86%  "

$priceCell = $ws.Range("D25")
$priceCell.NumberFormat = "@"
$priceCell.Value = "2.041"
$priceCell.Style = "Normal"
$ws.Range("E25").Value = "  +1.16%  "

$priceCell = $ws.Range("D26")
$priceCell.NumberFormat = "@"
$priceCell.Value = "1.999.33"
$priceCell.Style = "Normal"
$ws.Range("E26").Value = "  -5.59%  "

$ws.Range("E27").Value = "  -0.35%  "

$priceCell = $ws.Range("D28")
$priceCell.NumberFormat = "@"
$priceCell.Value = "18.30"
$priceCell.Style = "Normal"
$ws.Range("E28").Value = "  -0.85%  "

$priceCell = $ws.Range("D29")
$priceCell.NumberFormat = "@"
$priceCell.Value = "5.255"
$priceCell.Style = "Normal"
$ws.Range("E29").Value = "  -1.29%  "

$priceCell = $ws.Range("D30")
$priceCell.NumberFormat = "@"
$priceCell.Value = "1.836"
$priceCell.Style = "Normal"
$ws.Range("E30").Value = "  +6.23%  "

$priceCell = $ws.Range("D31")
$priceCell.NumberFormat = "@"
$priceCell.Value = "111.05"
$priceCell.Style = "Normal"
$ws.Range("E31").Value = "  -2.68%  "

$priceCell = $ws.Range("D32")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.08869"
$priceCell.Style = "Normal"
$ws.Range("E32").Value = "  -1.81%  "

$priceCell = $ws.Range("D33")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.7784"
$priceCell.Style = "Normal"
$ws.Range("E33").Value = "  -2.80%  "

$priceCell = $ws.Range("D34")
$priceCell.NumberFormat = "@"
$priceCell.Value = "4.565"
$priceCell.Style = "Normal"
$ws.Range("E34").Value = "  -5.81%  "

$priceCell = $ws.Range("D35")
$priceCell.NumberFormat = "@"
$priceCell.Value = "2.944"
$priceCell.Style = "Normal"
$ws.Range("E35").Value = "  -0.49%  "

$priceCell = $ws.Range("D36")
$priceCell.NumberFormat = "@"
$priceCell.Value = "1.095"
$priceCell.Style = "Normal"
$ws.Range("E36").Value = "  -6.46%  "

$priceCell = $ws.Range("D37")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.9987"
$priceCell.Style = "Normal"
$ws.Range("E37").Value = "  -0.27%  "

$priceCell = $ws.Range("D38")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.05438"
$priceCell.Style = "Normal"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("E39").Value = "  -2.09%  "

$priceCell = $ws.Range("D40")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.01929"
$priceCell.Style = "Normal"
$ws.Range("E40").Value = "  -1.26%  "

$priceCell = $ws.Range("D41")
$priceCell.NumberFormat = "@"
$priceCell.Value = "2.817"
$priceCell.Style = "Normal"
$ws.Range("E41").Value = "  -5.02%  "

$priceCell = $ws.Range("D42")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.5076"
$priceCell.Style = "Normal"
$ws.Range("E42").Value = "  -3.24%  "

$priceCell = $ws.Range("D43")
$priceCell.NumberFormat = "@"
$priceCell.Value = "6.815"
$priceCell.Style = "Normal"
$ws.Range("E43").Value = "  -4.45%  "

$priceCell = $ws.Range("D44")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.1642"
$priceCell.Style = "Normal"
$ws.Range("E44").Value = "  -1.87%  "

$priceCell = $ws.Range("D45")
$priceCell.NumberFormat = "@"
$priceCell.Value = "8.246"
$priceCell.Style = "Normal"
$ws.Range("E45").Value = "  -5.20%  "

$priceCell = $ws.Range("D46")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.06640"
$priceCell.Style = "Normal"
$ws.Range("E46").Value = "  -2.02%  "

$priceCell = $ws.Range("D47")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.4728"
$priceCell.Style = "Normal"
$ws.Range("E47").Value = "  -2.89%  "

$priceCell = $ws.Range("D48")
$priceCell.NumberFormat = "@"
$priceCell.Value = "105.93"
$priceCell.Style = "Normal"
$ws.Range("E48").Value = "  -1.73%  "

$priceCell = $ws.Range("D49")
$priceCell.NumberFormat = "@"
$priceCell.Value = "10.35"
$priceCell.Style = "Normal"
$ws.Range("E49").Value = "  -1.66%  "

$priceCell = $ws.Range("D50")
$priceCell.NumberFormat = "@"
$priceCell.Value = "0.9989"
$priceCell.Style = "Normal"
$ws.Range("E50").Value = "  -0.19%  "

$priceCell = $ws.Range("D51")
$priceCell.NumberFormat = "@"
$priceCell.Value = "1.643"
$priceCell.Style = "Normal"
$ws.Range("E51").Value = "  -2.16%  "
